$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (matching the original inline-string
# cells) without leaving a residual NumberFormat/style change behind.
# Excel auto-detects numeric-looking / percent-looking strings and
# converts them to Number cells unless the cell is formatted as Text
# first; ClearFormats() afterwards drops the now-unneeded "@" format so
# the cell's style index is left exactly as it was before the write.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "289.43"
Set-TextValue "E2" "-0.30%"
Set-TextValue "D3" "31.08"
Set-TextValue "E3" "1.12%"
Set-TextValue "D4" "4.943"
Set-TextValue "E4" "0.50%"
Set-TextValue "D5" "0.07352"
Set-TextValue "E5" "1.51%"
Set-TextValue "D6" "2.294"
Set-TextValue "E6" "16.87%"
Set-TextValue "E7" "0.36%"
Set-TextValue "D8" "0.9186"
Set-TextValue "E8" "1.97%"
Set-TextValue "D9" "0.09131"
Set-TextValue "E9" "13.88%"
Set-TextValue "E10" "1.24%"
Set-TextValue "D11" "0.08176"
Set-TextValue "E11" "1.11%"
Set-TextValue "D12" "0.03120"
Set-TextValue "E12" "2.05%"
Set-TextValue "D13" "0.09991"
Set-TextValue "E13" "-0.38%"
Set-TextValue "E14" "0.08%"
Set-TextValue "D15" "0.005757"
Set-TextValue "E15" "1.12%"
Set-TextValue "E16" "-0.10%"
Set-TextValue "D17" "3.745"
Set-TextValue "E17" "1.26%"
Set-TextValue "D19" "0.3330"
Set-TextValue "E19" "0.44%"
Set-TextValue "E20" "-0.29%"
Set-TextValue "D21" "4.189"
Set-TextValue "E21" "5.62%"
Set-TextValue "E22" "-2.09%"
Set-TextValue "D23" "0.04515"
Set-TextValue "E23" "-0.20%"
Set-TextValue "E24" "0.08%"
Set-TextValue "D25" "0.004197"
Set-TextValue "E25" "-5.10%"
Set-TextValue "E26" "0.07%"
Set-TextValue "E27" "0.02%"
Set-TextValue "D39" "0.01578"
Set-TextValue "E39" "-0.77%"
Set-TextValue "D40" "0.04502"
Set-TextValue "E40" "3.54%"
Set-TextValue "D41" "0.007377"
Set-TextValue "E41" "1.34%"
Set-TextValue "D42" "0.009857"
Set-TextValue "E42" "-1.83%"
Set-TextValue "D43" "0.1335"
Set-TextValue "E43" "1.58%"
Set-TextValue "E44" "10.79%"
Set-TextValue "D45" "0.008488"
Set-TextValue "E45" "-10.32%"
Set-TextValue "D46" "0.00006112"
Set-TextValue "E46" "5.24%"
Set-TextValue "E47" "0.01%"
Set-TextValue "D48" "2.445"
Set-TextValue "E48" "8.47%"
Set-TextValue "E49" "-30.99%"
Set-TextValue "E50" "0.01%"
Set-TextValue "E51" "0.01%"
